$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 593.2917
$ws.Range("I2").Value = 186.05556
$ws.Range("J2").Value = 1815
$ws.Range("K2").Value = 186.05556
$ws.Range("L2").Value = 1815
$ws.Range("M2").Value = -73.05556000000001
$ws.Range("N2").Value = -2041
$ws.Range("H18").Value = 1137.25
$ws.Range("I18").Value = 1137.25
$ws.Range("K18").Value = 1137.25
$ws.Range("M18").Value = -853.25
$ws.Range("H64").Value = 4936.625
$ws.Range("I64").Value = 4623.5
$ws.Range("K64").Value = 4623.5
$ws.Range("M64").Value = -4375.5
$ws.Range("H67").Value = 4936.625
$ws.Range("I67").Value = 4623.5
$ws.Range("K67").Value = 4623.5
$ws.Range("M67").Value = -3765.5
$ws.Range("H98").Value = 886.6923
$ws.Range("I98").Value = 958.8182
$ws.Range("J98").Value = 490
$ws.Range("K98").Value = 958.8182
$ws.Range("L98").Value = 490
$ws.Range("M98").Value = 539.1818
$ws.Range("N98").Value = -3486
$ws.Range("H122").Value = 886.6923
$ws.Range("I122").Value = 958.8182
$ws.Range("J122").Value = 490
$ws.Range("K122").Value = 2876.4546
$ws.Range("L122").Value = 1470
$ws.Range("M122").Value = -426.4546
$ws.Range("N122").Value = -6370
$ws.Range("H132").Value = 3798.7273
$ws.Range("I132").Value = 2420.889
$ws.Range("K132").Value = 7262.667
$ws.Range("M132").Value = -4732.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2296.7297
$ws.Range("I32").Value = 2417.8572
$ws.Range("K32").Value = 2417.8572
$ws.Range("M32").Value = -2130.8572
$ws.Range("H61").Value = 2802.2354
$ws.Range("I61").Value = 2802.2354
$ws.Range("K61").Value = 2802.2354
$ws.Range("M61").Value = -2590.2354
$ws.Range("H74").Value = 1924.5333
$ws.Range("I74").Value = 1764
$ws.Range("K74").Value = 1764
$ws.Range("M74").Value = -890
$ws.Range("H77").Value = 1924.5333
$ws.Range("I77").Value = 1764
$ws.Range("K77").Value = 8820
$ws.Range("M77").Value = -4452
$ws.Range("H136").Value = 2802.2354
$ws.Range("I136").Value = 2802.2354
$ws.Range("K136").Value = 8406.706200000001
$ws.Range("M136").Value = -5856.706200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 305.5
$ws.Range("I22").Value = 317.9091
$ws.Range("J22").Value = 260
$ws.Range("K22").Value = 317.9091
$ws.Range("L22").Value = 260
$ws.Range("M22").Value = -144.9091
$ws.Range("N22").Value = -606
$ws.Range("H105").Value = 4836
$ws.Range("I105").Value = 4698
$ws.Range("K105").Value = 4698
$ws.Range("M105").Value = -2951
$ws.Range("H134").Value = 3160.25
$ws.Range("I134").Value = 3116.0527
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 9348.158100000001
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -6813.158100000001
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3585.7273
$ws.Range("J105").Value = 3999.8333
$ws.Range("L105").Value = 3999.8333
$ws.Range("N105").Value = -7493.8333
$ws.Range("H132").Value = 2398.6191
$ws.Range("I132").Value = 2453.3845
$ws.Range("K132").Value = 7360.1535
$ws.Range("M132").Value = -4830.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 116402.42
$ws.Range("I2").Value = 100023.18
$ws.Range("J2").Value = 138923.88
$ws.Range("K2").Value = 600139.08
$ws.Range("L2").Value = 833543.28
$ws.Range("M2").Value = -600026.08
$ws.Range("N2").Value = -833769.28
$ws.Range("H68").Value = 8063340.5
$ws.Range("I68").Value = 494.66666
$ws.Range("K68").Value = 1483.99998
$ws.Range("M68").Value = -672.9999800000001
$ws.Range("H71").Value = 8063340.5
$ws.Range("I71").Value = 494.66666
$ws.Range("K71").Value = 4451.99994
$ws.Range("M71").Value = -395.9999399999997
$ws.Range("H88").Value = 10148.429
$ws.Range("J88").Value = 10508.2
$ws.Range("L88").Value = 31524.6
$ws.Range("N88").Value = -32380.6
$ws.Range("H91").Value = 10148.429
$ws.Range("J91").Value = 10508.2
$ws.Range("L91").Value = 31524.6
$ws.Range("N91").Value = -34488.60000000001
$ws.Range("H109").Value = 2750
$ws.Range("I109").Value = 2750
$ws.Range("K109").Value = 8250
$ws.Range("M109").Value = -7210
$ws.Range("H113").Value = 1472.7333
$ws.Range("J113").Value = 1392
$ws.Range("L113").Value = 4176
$ws.Range("N113").Value = -8516
$ws.Range("H129").Value = 837508.3
$ws.Range("I129").Value = 3559.25
$ws.Range("J129").Value = 1254482.9
$ws.Range("K129").Value = 10677.75
$ws.Range("L129").Value = 3763448.7
$ws.Range("M129").Value = -5677.75
$ws.Range("N129").Value = -3773448.7
$ws.Range("H131").Value = 835769.3
$ws.Range("J131").Value = 1002703.3
$ws.Range("L131").Value = 3008109.9
$ws.Range("N131").Value = -3018189.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 25040.428
$ws.Range("J98").Value = 25040.428
$ws.Range("L98").Value = 25040.428
$ws.Range("N98").Value = -31030.428
$ws.Range("H102").Value = 1727.25
$ws.Range("I102").Value = 1803.3334
$ws.Range("J102").Value = 1499
$ws.Range("K102").Value = 1803.3334
$ws.Range("L102").Value = 1499
$ws.Range("M102").Value = -181.3334
$ws.Range("N102").Value = -4743
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 1797.2858
$ws.Range("I132").Value = 1738.2
$ws.Range("J132").Value = 1945
$ws.Range("K132").Value = 5214.6
$ws.Range("L132").Value = 5835
$ws.Range("M132").Value = -2684.6
$ws.Range("N132").Value = -10895
$ws.Range("H133").Value = 59999
$ws.Range("J133").Value = 59999
$ws.Range("L133").Value = 59999
$ws.Range("N133").Value = -70119

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6024.227
$ws.Range("I7").Value = 2221.1428
$ws.Range("K7").Value = 2221.1428
$ws.Range("M7").Value = -2109.1428
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 4
$ws.Range("M14").Value = 168
$ws.Range("H16").Value = 1239.091
$ws.Range("I16").Value = 1141.5
$ws.Range("K16").Value = 1141.5
$ws.Range("M16").Value = -971.5
$ws.Range("H55").Value = 262.4
$ws.Range("I55").Value = 224
$ws.Range("K55").Value = 224
$ws.Range("M55").Value = -51
$ws.Range("H82").Value = 1320.6923
$ws.Range("I82").Value = 1580.6666
$ws.Range("J82").Value = 735.75
$ws.Range("K82").Value = 1580.6666
$ws.Range("L82").Value = 735.75
$ws.Range("M82").Value = -1219.6666
$ws.Range("N82").Value = -1457.75
$ws.Range("H85").Value = 1320.6923
$ws.Range("I85").Value = 1580.6666
$ws.Range("J85").Value = 735.75
$ws.Range("K85").Value = 1580.6666
$ws.Range("L85").Value = 735.75
$ws.Range("M85").Value = -332.6666
$ws.Range("N85").Value = -3231.75
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = 0
$ws.Range("H122").Value = 6048.857
$ws.Range("I122").Value = 6053.3184
$ws.Range("J122").Value = 6032.5
$ws.Range("K122").Value = 18159.9552
$ws.Range("L122").Value = 18097.5
$ws.Range("M122").Value = -15709.9552
$ws.Range("N122").Value = -22997.5
$ws.Range("H126").Value = 6024.227
$ws.Range("I126").Value = 2221.1428
$ws.Range("K126").Value = 6663.428400000001
$ws.Range("M126").Value = -4193.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 16600
$ws.Range("J18").Value = 16600
$ws.Range("L18").Value = 16600
$ws.Range("N18").Value = -16946
$ws.Range("H51").Value = 28750
$ws.Range("I51").Value = 25000
$ws.Range("J51").Value = 31000
$ws.Range("K51").Value = 25000
$ws.Range("L51").Value = 31000
$ws.Range("M51").Value = -24490
$ws.Range("N51").Value = -32020
$ws.Range("H52").Value = 18999.666
$ws.Range("I52").Value = 5000
$ws.Range("J52").Value = 21799.6
$ws.Range("K52").Value = 5000
$ws.Range("L52").Value = 21799.6
$ws.Range("M52").Value = -4774
$ws.Range("N52").Value = -22251.6
$ws.Range("H62").Value = 18679.4
$ws.Range("I62").Value = 18933.666
$ws.Range("J62").Value = 18570.428
$ws.Range("K62").Value = 18933.666
$ws.Range("L62").Value = 18570.428
$ws.Range("M62").Value = -18309.666
$ws.Range("N62").Value = -19818.428
$ws.Range("H65").Value = 18679.4
$ws.Range("I65").Value = 18933.666
$ws.Range("J65").Value = 18570.428
$ws.Range("K65").Value = 94668.33
$ws.Range("L65").Value = 92852.14
$ws.Range("M65").Value = -91548.33
$ws.Range("N65").Value = -99092.14
$ws.Range("H132").Value = 5872.706
$ws.Range("I132").Value = 5655.933
$ws.Range("K132").Value = 16967.799
$ws.Range("M132").Value = -14437.799
$ws.Range("H136").Value = 4021.2632
$ws.Range("I136").Value = 4021.2632
$ws.Range("K136").Value = 12063.7896
$ws.Range("M136").Value = -9513.7896
